$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Move existing STOK columns (H1:J1 -> Q1:S1), preserving value + style/format
$ws.Range("H1:J1").Copy($ws.Range("Q1"))

# 2. Set the new header text, in the order needed so the shared-string table
#    is built up in the same sequence as the target workbook.
$ws.Range("H1").Value = "HARGA GROSIR #1"
$ws.Range("K1").Value = "HARGA GROSIR #2"
$ws.Range("N1").Value = "HARGA GROSIR #3"
$ws.Range("I1").Value = "JUMLAH MINIMAL GROSIR #1"
$ws.Range("L1").Value = "JUMLAH MINIMAL GROSIR #2"
$ws.Range("J1").Value = "DISKON GROSIR #1"
$ws.Range("M1").Value = "DISKON GROSIR #2"
$ws.Range("O1").Value = "JUMLAH MINIMAL GROSIR #3"
$ws.Range("P1").Value = "DISKON GROSIR #3"
$ws.Range("F1").Value = "HARGA JUAL (RETAIL)"
$ws.Range("G1").Value = "DISKON (RETAIL)"

# 3. Apply the header format (same as the other header cells) to the brand
#    new cells K1:P1.
$ws.Range("H1").Copy()
$ws.Range("K1:P1").PasteSpecial(-4122)

# 4. Resize the columns to match the new layout. The ColumnWidth property
#    is internally quantized by the engine to steps of 1/6 of a character,
#    offset by 5/6 (the standard 5-pixel cell padding) - back out the
#    character width that reproduces each target stored sheet width.
function ColWidthFor($targetStoredWidth) {
    $n = [Math]::Round(($targetStoredWidth - 5.0/6.0) * 6.0)
    return $n / 6.0
}

$ws.Columns.Item(6).ColumnWidth = ColWidthFor(19.85546875)
$ws.Columns.Item(7).ColumnWidth = ColWidthFor(15.7109375)
$ws.Columns.Item(8).ColumnWidth = ColWidthFor(16.85546875)
$ws.Columns.Item(9).ColumnWidth = ColWidthFor(26.7109375)
$ws.Columns.Item(10).ColumnWidth = ColWidthFor(20.140625)
$ws.Columns.Item(11).ColumnWidth = ColWidthFor(16.85546875)
$ws.Columns.Item(12).ColumnWidth = ColWidthFor(26.7109375)
$ws.Columns.Item(13).ColumnWidth = ColWidthFor(17.5703125)
$ws.Columns.Item(14).ColumnWidth = ColWidthFor(16.85546875)
$ws.Columns.Item(15).ColumnWidth = ColWidthFor(26.7109375)
$ws.Columns.Item(16).ColumnWidth = ColWidthFor(17.5703125)
$ws.Columns.Item(17).ColumnWidth = ColWidthFor(13.5703125)
$ws.Columns.Item(18).ColumnWidth = ColWidthFor(14)
$ws.Columns.Item(19).ColumnWidth = ColWidthFor(22.85546875)

# 5. Update the active selection shown in the sheet view.
$ws.Range("E5").Select()
